$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.281.57"
$ws.Range("E2").Value = "  -5.66%  "
$ws.Range("D3").Value = "2.223.70"
$ws.Range("E3").Value = "  -5.34%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.12"
$ws.Range("E5").Value = "  +2.39%  "
$ws.Range("E6").Value = "  -6.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "70.30"
$ws.Range("E7").Value = "  -5.40%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -7.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.96"
$ws.Range("E10").Value = "  +4.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0952"
$ws.Range("E11").Value = "  -6.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.42"
$ws.Range("E12").Value = "  -1.76%  "
$ws.Range("E13").Value = "  -3.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.77"
$ws.Range("E14").Value = "  -7.63%  "
$ws.Range("D15").Value = "2.551.22"
$ws.Range("E15").Value = "  -5.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.82"
$ws.Range("E16").Value = "  -9.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.844"
$ws.Range("E17").Value = "  -8.68%  "
$ws.Range("D18").Value = "2.217.48"
$ws.Range("E18").Value = "  -5.87%  "
$ws.Range("D19").Value = "41.241.78"
$ws.Range("E19").Value = "  -5.56%  "
$ws.Range("D20").Value = [string]::Concat("0.0", [char]0x2083, "0953")
$ws.Range("E20").Value = "  -8.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.50"
$ws.Range("E21").Value = "  -5.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.09"
$ws.Range("E22").Value = "  -7.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.22"
$ws.Range("E23").Value = "  -8.26%  "
$ws.Range("E24").Value = "  +11.61%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  -2.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.42"
$ws.Range("E27").Value = "  -2.70%  "
$ws.Range("E28").Value = "  -7.17%  "
$ws.Range("E29").Value = "  -4.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.04"
$ws.Range("E30").Value = "  -1.71%  "
$ws.Range("E31").Value = "  -7.82%  "
$ws.Range("E32").Value = "  -7.34%  "
$ws.Range("E33").Value = "  -6.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0713"
$ws.Range("E34").Value = "  -5.51%  "
$ws.Range("E35").Value = "  -4.61%  "
$ws.Range("E36").Value = "  -9.98%  "
$ws.Range("E37").Value = "  +3.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.31"
$ws.Range("E38").Value = "  +16.35%  "
$ws.Range("E39").Value = "  -1.26%  "
$ws.Range("E40").Value = "  -5.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.85"
$ws.Range("E41").Value = "  -11.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.45"
$ws.Range("E42").Value = "  -1.41%  "
$ws.Range("E43").Value = "  -8.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.204"
$ws.Range("E44").Value = "  +1.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.81"
$ws.Range("E45").Value = "  -2.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.94"
$ws.Range("E46").Value = "  +10.98%  "
$ws.Range("E47").Value = "  -6.61%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("E49").Value = "  +4.87%  "
$ws.Range("E50").Value = "  -5.88%  "
$ws.Range("E51").Value = "  -4.93%  "
